$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K16: long list of raw ICD codes used by the new "unique icd" filtering check.
# The stored text itself starts with "465.8" (no leading quote) but the cell is
# left in Excel's "quote prefix" (force-text) state with wrap-left alignment.
# We reproduce that by typing a leading apostrophe -- Excel strips it from the
# stored value but keeps the quotePrefix flag on the cell's style -- and then
# turning wrapping on.
$icdLong = @"
'465.8',
 '465.9',
 '466',
 '491.21',
 '493',
 '493.01',
 '493.02',
 '493.1',
 '493.11',
 '493.2',
 '493.81',
 '493.82',
 '493.9',
 '493.91',
 '493.92',
 'J06.9',
 'J20.9',
 'J44.1',
 'J44.9',
 'J45.20',
 'J45.21',
 'J45.22',
 'J45.30',
 'J45.31',
 'J45.32',
 'J45.40',
 'J45.41',
 'J45.42',
 'J45.50',
 'J45.51',
 'J45.52',
 'J45.901',
 'J45.902',
 'J45.909',
 'J45.991',
 'J45.998'
"@
$ws.Range("K16").Value = $icdLong
$ws.Range("K16").WrapText = $true

# Typing the long, wrapped ICD text auto-expands row 16's height; the row was
# already a fixed (customHeight) 120.75pt row before this edit and stays that
# way afterward, so put it back.
$ws.Rows.Item(16).RowHeight = 120.75

# New header columns (row 1): L = import_unique_icd, M = check_unique_icd.
$ws.Range("L1").Value = "import_unique_icd"
$ws.Range("M1").Value = "check_unique_icd"

# Remaining new numeric check/import counts on row 16.
$ws.Range("G16").Value = 13241
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = 3577
$ws.Range("J16").Value = 75

# M16: short, de-duplicated ICD list for the new "unique icd" check column.
# Font is switched to Arial Unicode MS 10pt (distinct from the sheet default),
# and the cell keeps General alignment (reset via the Normal style) rather than
# the column's inherited left alignment. The text itself keeps its literal
# leading apostrophe, so it is typed doubled (Excel's force-text trigger
# character plus the literal character we want to keep).
$icdShort = "''466', '493', '493.1', '493.2', '493.9'"
$ws.Range("M16").Value = $icdShort
$ws.Range("M16").Style = "Normal"
$ws.Range("M16").Font.Size = 10
$ws.Range("M16").Font.Name = "Arial Unicode MS"

# Final cursor/selection position left by the editing session.
$ws.Range("C5").Select()
